$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G1").Value = "VenueID"
$ws.Range("G2").Value = 1
$ws.Range("G3").Select()
